# Auto-generated edit script applying cell value changes per diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 48500
$ws.Range("J3").Value = 48500
$ws.Range("L3").Value = 48500
$ws.Range("N3").Value = -48728

$ws.Range("H4").Value = 362.5
$ws.Range("I4").Value = 150
$ws.Range("J4").Value = 1000
$ws.Range("K4").Value = 150
$ws.Range("L4").Value = 1000
$ws.Range("M4").Value = -36
$ws.Range("N4").Value = -1228

$ws.Range("H5").Value = 100.84615
$ws.Range("I5").Value = 50.2
$ws.Range("J5").Value = 132.5
$ws.Range("K5").Value = 50.2
$ws.Range("L5").Value = 132.5
$ws.Range("M5").Value = 64.8
$ws.Range("N5").Value = -362.5

$ws.Range("H11").Value = 72599.21000000001
$ws.Range("I11").Value = 72599.21000000001
$ws.Range("K11").Value = 72599.21000000001
$ws.Range("M11").Value = -72459.21000000001

$ws.Range("H28").Value = 460.68182
$ws.Range("I28").Value = 329.05884
$ws.Range("J28").Value = 908.2
$ws.Range("K28").Value = 329.05884
$ws.Range("L28").Value = 908.2
$ws.Range("M28").Value = 155.94116
$ws.Range("N28").Value = -1878.2

$ws.Range("H41").Value = 1313.8182
$ws.Range("I41").Value = 300
$ws.Range("J41").Value = 1415.2
$ws.Range("K41").Value = 300
$ws.Range("L41").Value = 1415.2
$ws.Range("M41").Value = 140
$ws.Range("N41").Value = -2295.2

$ws.Range("H102").Value = 48500
$ws.Range("J102").Value = 48500
$ws.Range("L102").Value = 48500
$ws.Range("N102").Value = -54990

$ws.Range("H137").Value = 1498.2693
$ws.Range("I137").Value = 1970.4073
$ws.Range("K137").Value = 5911.2219
$ws.Range("M137").Value = -3361.2219

$ws.Range("H138").Value = 1220.2549
$ws.Range("I138").Value = 1034.3334
$ws.Range("J138").Value = 1561.1111
$ws.Range("K138").Value = 3103.0002
$ws.Range("L138").Value = 4683.3333
$ws.Range("M138").Value = 2036.9998
$ws.Range("N138").Value = -14963.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 92.5
$ws.Range("I5").Value = 85
$ws.Range("K5").Value = 85
$ws.Range("M5").Value = 27

$ws.Range("H32").Value = 4759.34
$ws.Range("I32").Value = 4425.05
$ws.Range("J32").Value = 6096.5
$ws.Range("K32").Value = 4425.05
$ws.Range("L32").Value = 6096.5
$ws.Range("M32").Value = -4138.05
$ws.Range("N32").Value = -6670.5

$ws.Range("H132").Value = 4163.6333
$ws.Range("I132").Value = 3523.318
$ws.Range("J132").Value = 4534.3423
$ws.Range("K132").Value = 10569.954
$ws.Range("L132").Value = 13603.0269
$ws.Range("M132").Value = -8039.954000000002
$ws.Range("N132").Value = -18663.0269

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 92.5
$ws.Range("I4").Value = 85
$ws.Range("K4").Value = 85
$ws.Range("M4").Value = 30

$ws.Range("H22").Value = 2018.2273
$ws.Range("I22").Value = 1531.6316
$ws.Range("J22").Value = 5100
$ws.Range("K22").Value = 1531.6316
$ws.Range("L22").Value = 5100
$ws.Range("M22").Value = -1358.6316
$ws.Range("N22").Value = -5446

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 18501500
$ws.Range("I6").Value = 22201200
$ws.Range("J6").Value = 3000
$ws.Range("K6").Value = 22201200
$ws.Range("L6").Value = 3000
$ws.Range("M6").Value = -22201087
$ws.Range("N6").Value = -3226

$ws.Range("H23").Value = 24782.912
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 24782.912
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 24782.912
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -25262.912

$ws.Range("H27").Value = 24782.912
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 24782.912
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 24782.912
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -25166.912

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 2116.2727
$ws.Range("I13").Value = 426.33334
$ws.Range("J13").Value = 2750
$ws.Range("K13").Value = 1279.00002
$ws.Range("L13").Value = 8250
$ws.Range("M13").Value = -1111.00002
$ws.Range("N13").Value = -8586

$ws.Range("H23").Value = 4947.952
$ws.Range("I23").Value = 83.333336
$ws.Range("J23").Value = 5758.722
$ws.Range("K23").Value = 250.000008
$ws.Range("L23").Value = 17276.166
$ws.Range("M23").Value = -15.00000800000001
$ws.Range("N23").Value = -17746.166

$ws.Range("H68").Value = 415.83334
$ws.Range("I68").Value = 301.84616
$ws.Range("J68").Value = 550.5454999999999
$ws.Range("K68").Value = 905.5384799999999
$ws.Range("L68").Value = 1651.6365
$ws.Range("M68").Value = -94.53847999999994
$ws.Range("N68").Value = -3273.6365

$ws.Range("H71").Value = 415.83334
$ws.Range("I71").Value = 301.84616
$ws.Range("J71").Value = 550.5454999999999
$ws.Range("K71").Value = 2716.61544
$ws.Range("L71").Value = 4954.9095
$ws.Range("M71").Value = 1339.38456
$ws.Range("N71").Value = -13066.9095

$ws.Range("H100").Value = 5028
$ws.Range("J100").Value = 5028
$ws.Range("L100").Value = 15084
$ws.Range("N100").Value = -16706

$ws.Range("H106").Value = 4115
$ws.Range("J106").Value = 4115
$ws.Range("L106").Value = 12345
$ws.Range("N106").Value = -14237

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2501498
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 2501498
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 2501498
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -2501874

$ws.Range("H55").Value = 262.6875
$ws.Range("I55").Value = 180.08333
$ws.Range("J55").Value = 510.5
$ws.Range("K55").Value = 180.08333
$ws.Range("L55").Value = 510.5
$ws.Range("M55").Value = -7.083329999999989
$ws.Range("N55").Value = -856.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 52082.125
$ws.Range("J46").Value = 52082.125
$ws.Range("L46").Value = 52082.125
$ws.Range("N46").Value = -52544.125

$ws.Range("H134").Value = 52082.125
$ws.Range("J134").Value = 52082.125
$ws.Range("L134").Value = 156246.375
$ws.Range("N134").Value = -161316.375
